$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free approach: for "Price" (column D) cells whose new text looks like a
# number (e.g. "582.67"), Excel's COM layer auto-converts the assigned string to a
# numeric cell. To keep these cells as text (matching the original inlineStr/text
# cells) we temporarily force a text number-format, assign the value, then reset
# the cell style back to Normal so no stray formatting is left behind.

# --- Row 2: Bitcoin ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "72.292.65"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.32%  "

# --- Row 3: Ethereum ---
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.962.28"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.27%  "

# --- Row 4: TetherUSD ---
$ws.Range("E4").Value = "  +0.11%  "

# --- Row 5: BNB ---
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +10.84%  "

# --- Row 6: Solana ---
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.89%  "

# --- Row 7: XRP ---
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.677"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.94%  "

# --- Row 8: USDC ---
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.03%  "

# --- Row 9: Cardano ---
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.748"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.28%  "

# --- Row 10: Dogecoin ---
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.167"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.15%  "

# --- Row 11: Avalanche ---
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.02"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.06%  "

# --- Row 12: ShibaInu ---
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000317"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.96%  "

# --- Row 13: Polkadot ---
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.82"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.00%  "

# --- Row 14: WrappedliquidstakedEther2.0 ---
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.608.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.92%  "

# --- Row 15: WrappedEther ---
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.992.55"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.68%  "

# --- Row 16: Polygon ---
$ws.Range("E16").Value = "  +6.57%  "

# --- Row 17: Uniswap ---
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.94"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.33%  "

# --- Row 18: Chainlink ---
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.40"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.67%  "

# --- Row 19: TRON ---
$ws.Range("E19").Value = "  -0.70%  "

# --- Row 20: WrappedBTC ---
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.338.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.64%  "

# --- Row 21: BitcoinCash ---
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "432.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.68%  "

# --- Row 22: PancakeSwap ---
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +12.98%  "

# --- Row 23: Litecoin ---
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "95.62"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.64%  "

# --- Row 24: ImmutableX ---
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.80%  "

# --- Row 25: InternetComputer(DFINITY) ---
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.26%  "

# (Row 26: Toncoin -- unchanged)

# --- Row 27: RenderToken ---
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.58"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.43%  "

# --- Row 28: Filecoin ---
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.98%  "

# --- Row 29: LEO ---
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.92"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.96%  "

# --- Row 30: EthereumClassic ---
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.98%  "

# --- Row 31: NEARProtocol ---
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.96"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.21%  "

# --- Row 32: InjectiveProtocol ---
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.57"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.65%  "

# --- Row 33 / 34: Hedera <-> Cosmos swap position ---
$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "13.49"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.75%  "

$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.132"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.62%  "

# --- Row 35: Bittensor ---
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "680.31"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.15%  "

# --- Row 36: OKB ---
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "68.44"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.68%  "

# --- Row 37: TheGraph ---
$ws.Range("E37").Value = "  -1.33%  "

# --- Row 38: PEPE ---
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0852"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.70%  "

# --- Row 39 / 40: Kaspa <-> WEMIXToken swap position ---
$ws.Range("B39").Value = "WEMIXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.38"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.44%  "

$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.146"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.87%  "

# --- Row 41: Dai ---
$ws.Range("E41").Value = "  +0.00%  "

# --- Row 42: THORChain ---
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.06"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.89%  "

# --- Row 43 / 44: FirstDigitalUSD <-> ThetaToken swap position ---
$ws.Range("B43").Value = "ThetaToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.30"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.35%  "

$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.05%  "

# --- Row 45: VeChain ---
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0485"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.48%  "

# --- Row 46: Fetch.AI ---
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.77"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.61%  "

# --- Row 47: Stellar ---
$ws.Range("E47").Value = "  -0.19%  "

# --- Row 48: ApeXProtocol ---
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.35"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.62%  "

# --- Row 49: LidoDAOToken ---
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.45"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.02%  "

# --- Row 50: Stacks ---
$ws.Range("E50").Value = "  -0.19%  "

# --- Row 51: Maker ---
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.779.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +8.96%  "
